$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 15217.857
$ws.Range("I33").Value = 21083.4
$ws.Range("J33").Value = 554
$ws.Range("K33").Value = 21083.4
$ws.Range("L33").Value = 554
$ws.Range("M33").Value = -20854.4
$ws.Range("N33").Value = -1012

# Row 45
$ws.Range("H45").Value = 7000
$ws.Range("I45").Value = 7000
$ws.Range("J45").Value = 7000
$ws.Range("K45").Value = 21000
$ws.Range("L45").Value = 21000
$ws.Range("M45").Value = -20808
$ws.Range("N45").Value = -21384

# Row 98
$ws.Range("H98").Value = 1411.0588
$ws.Range("I98").Value = 1411.0588
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1411.0588
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 86.94119999999998
$ws.Range("N98").ClearContents()

# Row 122
$ws.Range("H122").Value = 1411.0588
$ws.Range("I122").Value = 1411.0588
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4233.1764
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1783.1764
$ws.Range("N122").ClearContents()

# Row 132
$ws.Range("H132").Value = 3030.6875
$ws.Range("I132").Value = 3030.6875
$ws.Range("K132").Value = 9092.0625
$ws.Range("M132").Value = -6562.0625

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9330.315000000001
$ws.Range("I32").Value = 6642.2764
$ws.Range("K32").Value = 6642.2764
$ws.Range("M32").Value = -6355.2764

$ws = $wb.Worksheets.Item("BSM")
# Row 103
$ws.Range("H103").Value = 15599.5
$ws.Range("J103").Value = 15599.5
$ws.Range("L103").Value = 15599.5
$ws.Range("N103").Value = -17943.5

# Row 106
$ws.Range("H106").Value = 671671
$ws.Range("J106").Value = 671671
$ws.Range("L106").Value = 671671
$ws.Range("N106").Value = -674195

$ws = $wb.Worksheets.Item("CRP")
# Row 21
$ws.Range("H21").Value = 13
$ws.Range("I21").Value = 13
$ws.Range("K21").Value = 13
$ws.Range("M21").Value = 222

# Row 69
$ws.Range("H69").Value = 58280.168
$ws.Range("I69").Value = 76920.25
$ws.Range("J69").Value = 21000
$ws.Range("K69").Value = 76920.25
$ws.Range("L69").Value = 21000
$ws.Range("M69").Value = -76171.25
$ws.Range("N69").Value = -22498

# Row 72
$ws.Range("H72").Value = 58280.168
$ws.Range("I72").Value = 76920.25
$ws.Range("J72").Value = 21000
$ws.Range("K72").Value = 230760.75
$ws.Range("L72").Value = 63000
$ws.Range("M72").Value = -227016.75
$ws.Range("N72").Value = -70488

# Row 99
$ws.Range("H99").Value = 13492.196
$ws.Range("I99").Value = 17805.1
$ws.Range("K99").Value = 17805.1
$ws.Range("M99").Value = -16307.1

# Row 105
$ws.Range("H105").Value = 11266.091
$ws.Range("I105").Value = 1276
$ws.Range("K105").Value = 1276
$ws.Range("M105").Value = 471

# Row 126
$ws.Range("H126").Value = 13492.196
$ws.Range("I126").Value = 17805.1
$ws.Range("K126").Value = 53415.3
$ws.Range("M126").Value = -50945.3

# Row 132
$ws.Range("H132").Value = 46995.668
$ws.Range("I132").Value = 62518.727
$ws.Range("J132").Value = 4307.25
$ws.Range("K132").Value = 187556.181
$ws.Range("L132").Value = 12921.75
$ws.Range("M132").Value = -185026.181
$ws.Range("N132").Value = -17981.75

# Row 141
$ws.Range("H141").Value = 235998.8
$ws.Range("J141").Value = 235998.8
$ws.Range("L141").Value = 235998.8
$ws.Range("N141").Value = -246358.8

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 25670538
$ws.Range("I4").Value = 32271726
$ws.Range("J4").Value = 16368866
$ws.Range("K4").Value = 96815178
$ws.Range("L4").Value = 49106598
$ws.Range("M4").Value = -96815066
$ws.Range("N4").Value = -49106822

# Row 110
$ws.Range("H110").Value = 4783.5
$ws.Range("I110").Value = 4567
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 13701
$ws.Range("L110").Value = 15000
$ws.Range("M110").Value = -9611
$ws.Range("N110").Value = -23180

# Row 128
$ws.Range("H128").Value = 115583.5
$ws.Range("I128").Value = 115583.5
$ws.Range("K128").Value = 346750.5
$ws.Range("M128").Value = -341770.5

# Row 131
$ws.Range("H131").Value = 1735.4828
$ws.Range("J131").Value = 1745.3462
$ws.Range("L131").Value = 5236.0386
$ws.Range("N131").Value = -15316.0386

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2687.4736
$ws.Range("I132").Value = 2566.375
$ws.Range("K132").Value = 7699.125
$ws.Range("M132").Value = -5169.125

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2646.0557
$ws.Range("I7").Value = 2710.4546
$ws.Range("K7").Value = 2710.4546
$ws.Range("M7").Value = -2598.4546

# Row 68
$ws.Range("H68").Value = 2800.2
$ws.Range("I68").Value = 2000.3334
$ws.Range("K68").Value = 2000.3334
$ws.Range("M68").Value = -1251.3334

# Row 71
$ws.Range("H71").Value = 2800.2
$ws.Range("I71").Value = 2000.3334
$ws.Range("K71").Value = 10001.667
$ws.Range("M71").Value = -6257.666999999999

# Row 122
$ws.Range("H122").Value = 4132.816
$ws.Range("I122").Value = 2745.7856
$ws.Range("J122").Value = 4941.9165
$ws.Range("K122").Value = 8237.356800000001
$ws.Range("L122").Value = 14825.7495
$ws.Range("M122").Value = -5787.356800000001
$ws.Range("N122").Value = -19725.7495

# Row 126
$ws.Range("H126").Value = 2646.0557
$ws.Range("I126").Value = 2710.4546
$ws.Range("K126").Value = 8131.3638
$ws.Range("M126").Value = -5661.3638

# Row 132
$ws.Range("H132").Value = 20711.207
$ws.Range("I132").Value = 23886.908
$ws.Range("K132").Value = 71660.724
$ws.Range("M132").Value = -69130.724

# Row 136
$ws.Range("H136").Value = 1844.0385
$ws.Range("I136").Value = 555.8333
$ws.Range("K136").Value = 1667.4999
$ws.Range("M136").Value = 882.5001

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 256666.67
$ws.Range("J2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("N2").Value = -10224

# Row 51
$ws.Range("H51").Value = 27306.615
$ws.Range("I51").Value = 20748.75
$ws.Range("K51").Value = 20748.75
$ws.Range("M51").Value = -20238.75

# Row 113
$ws.Range("H113").Value = 630.6111
$ws.Range("J113").Value = 1073.4286
$ws.Range("L113").Value = 3220.2858
$ws.Range("N113").Value = -7560.2858

# Row 126
$ws.Range("H126").Value = 3690.65
$ws.Range("I126").Value = 4032.4443
$ws.Range("K126").Value = 12097.3329
$ws.Range("M126").Value = -9627.332900000001

# Row 136
$ws.Range("H136").Value = 4521.3335
$ws.Range("I136").Value = 1323.875
$ws.Range("K136").Value = 3971.625
$ws.Range("M136").Value = -1421.625

# Row 137
$ws.Range("H137").Value = 131989
$ws.Range("J137").Value = 131989
$ws.Range("L137").Value = 131989
$ws.Range("N137").Value = -142189
